# Update cryptos list with latest prices/volume figures (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.577.85"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.912.43"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.68"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.67"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.908.76"
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.23"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.568.24"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.928.27"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.718.90"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.45"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.16"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "485.65"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000169"
$ws.Range("E24").Value = "  +10.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.76"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.02"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.069.28"
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.80"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.02"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.862.17"
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.04"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.138"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "428.07"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.47"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.50"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.17"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.19"
$ws.Range("E49").Value = "  +6.30%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.823.12"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.27"
$ws.Range("E51").Value = "  -1.33%  "

Write-Output "Applied 96 cell updates"
